$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right before row 160, pushing the existing
# rows 160-226 down to 163-229 (new week of data is prepended at the
# top of this variety/quality block, most-recent-first ordering).
$ws.Rows.Item(160).Insert()
$ws.Rows.Item(160).Insert()
$ws.Rows.Item(160).Insert()

# Fill in the 3 new rows (160-162) with the new weekly report
# (fecha = 44466) for Especial / Primera / Segunda qualities.

# Row 160 - Especial
$ws.Cells.Item(160, 1).Value = 8
$ws.Cells.Item(160, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44466
$ws.Cells.Item(160, 5).Value = 4
$ws.Cells.Item(160, 6).Value = "Fruta"
$ws.Cells.Item(160, 7).Value = 100101
$ws.Cells.Item(160, 8).Value = "Berries"
$ws.Cells.Item(160, 9).Value = 100101007
$ws.Cells.Item(160, 10).Value = "Kiwi"
$ws.Cells.Item(160, 11).Value = "Hayward"
$ws.Cells.Item(160, 12).Value = "Especial"
$ws.Cells.Item(160, 13).Value = 16
$ws.Cells.Item(160, 14).Value = 410000
$ws.Cells.Item(160, 15).Value = 420000
$ws.Cells.Item(160, 16).Value = 415000
$ws.Cells.Item(160, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(160, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(160, 19).Value = 922
$ws.Cells.Item(160, 20).Value = 450

# Row 161 - Primera
$ws.Cells.Item(161, 1).Value = 8
$ws.Cells.Item(161, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(161, 3).Value = "Coquimbo"
$ws.Cells.Item(161, 4).Value = 44466
$ws.Cells.Item(161, 5).Value = 4
$ws.Cells.Item(161, 6).Value = "Fruta"
$ws.Cells.Item(161, 7).Value = 100101
$ws.Cells.Item(161, 8).Value = "Berries"
$ws.Cells.Item(161, 9).Value = 100101007
$ws.Cells.Item(161, 10).Value = "Kiwi"
$ws.Cells.Item(161, 11).Value = "Hayward"
$ws.Cells.Item(161, 12).Value = "Primera"
$ws.Cells.Item(161, 13).Value = 20
$ws.Cells.Item(161, 14).Value = 310000
$ws.Cells.Item(161, 15).Value = 320000
$ws.Cells.Item(161, 16).Value = 315000
$ws.Cells.Item(161, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(161, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(161, 19).Value = 700
$ws.Cells.Item(161, 20).Value = 450

# Row 162 - Segunda
$ws.Cells.Item(162, 1).Value = 8
$ws.Cells.Item(162, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(162, 3).Value = "Coquimbo"
$ws.Cells.Item(162, 4).Value = 44466
$ws.Cells.Item(162, 5).Value = 4
$ws.Cells.Item(162, 6).Value = "Fruta"
$ws.Cells.Item(162, 7).Value = 100101
$ws.Cells.Item(162, 8).Value = "Berries"
$ws.Cells.Item(162, 9).Value = 100101007
$ws.Cells.Item(162, 10).Value = "Kiwi"
$ws.Cells.Item(162, 11).Value = "Hayward"
$ws.Cells.Item(162, 12).Value = "Segunda"
$ws.Cells.Item(162, 13).Value = 20
$ws.Cells.Item(162, 14).Value = 270000
$ws.Cells.Item(162, 15).Value = 280000
$ws.Cells.Item(162, 16).Value = 275000
$ws.Cells.Item(162, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(162, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(162, 19).Value = 611
$ws.Cells.Item(162, 20).Value = 450
